$wb = $excel.ActiveWorkbook

$oldGuidFile = "ea1f9ec7-2a1c-4e23-8c73-e8cdc0c1f43e.md"
$newGuidFile = "f639bba2-5216-475f-bfc3-d5bc0d3b96fd.md"

$oldZhXlf = "ea1f9ec7-2a1c-4e23-8c73-e8cdc0c1f43e.ac0968c45c70563e6dcd6785ce3b86a66ba5b881.zh-cn.xlf"
$newZhXlf = "f639bba2-5216-475f-bfc3-d5bc0d3b96fd.2fa5403ca8cd31adebfc0cb914365cd909838f74.zh-cn.xlf"

$oldDeXlf = "ea1f9ec7-2a1c-4e23-8c73-e8cdc0c1f43e.ac0968c45c70563e6dcd6785ce3b86a66ba5b881.de-de.xlf"
$newDeXlf = "f639bba2-5216-475f-bfc3-d5bc0d3b96fd.2fa5403ca8cd31adebfc0cb914365cd909838f74.de-de.xlf"

$newLatestHandoffDate = "2016-03-24 11:57:01"
$newZhHandoffDatetime = "2016-03-24 11:56:57"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuidFile
$wsOverview.Range("D2").Value = $newLatestHandoffDate
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newGuidFile
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newGuidFile
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhHandoffDatetime
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldGuidFile) {
        $hl.TextToDisplay = $newGuidFile
    } else {
        $hl.TextToDisplay = $newZhXlf
    }
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newGuidFile
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newLatestHandoffDate
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldGuidFile) {
        $hl.TextToDisplay = $newGuidFile
    } else {
        $hl.TextToDisplay = $newDeXlf
    }
}
